# Natmi following Dr Hou advice
# Update Sema4d-Met LR-pair sheet: recompute edge stats for the full 3x3
# Sending-cluster x Target-cluster grid (ECs / FAPs / sCs), replacing the old
# 6-row partial grid with the complete 9-row grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("ECs", "Sema4d", "Met", "ECs", 3, 1, 68.50681433333334, 205.520443, 0.9663865053086182, 0.9663865053086185, 3, 1, 2.254050666666667, 6.762152, 0.5759961839619929, 0.575996183961993, 154.4178305192596, 1389.760474673336, 0.5566349392901303, 0.5566349392901305)
    3 = @("ECs", "Sema4d", "Met", "FAPs", 3, 1, 68.50681433333334, 205.520443, 0.9663865053086182, 0.9663865053086185, 1, 0.3333333333333333, 0.05371366666666667, 0.161141, 0.01372589688605336, 0.01372589688605336, 3.679752189495889, 33.117769705463, 0.01326452152393955, 0.01326452152393956)
    4 = @("ECs", "Sema4d", "Met", "sCs", 3, 1, 68.50681433333334, 205.520443, 0.9663865053086182, 0.9663865053086185, 3, 1, 1.605544, 4.816632, 0.4102779191519537, 0.4102779191519537, 109.9907047119974, 989.916342407976, 0.3964870444945483, 0.3964870444945485)
    5 = @("FAPs", "Sema4d", "Met", "ECs", 3, 1, 1.228643, 3.685929, 0.01733176511849816, 0.01733176511849817, 3, 1, 2.254050666666667, 6.762152, 0.5759961839619929, 0.575996183961993, 2.769423573245333, 24.924812159208, 0.009983030569580519, 0.009983030569580524)
    6 = @("FAPs", "Sema4d", "Met", "FAPs", 3, 1, 1.228643, 3.685929, 0.01733176511849816, 0.01733176511849817, 1, 0.3333333333333333, 0.05371366666666667, 0.161141, 0.01372589688605336, 0.01372589688605336, 0.06599492055433333, 0.5939542849889999, 0.0002378940208698022, 0.0002378940208698022)
    7 = @("FAPs", "Sema4d", "Met", "sCs", 3, 1, 1.228643, 3.685929, 0.01733176511849816, 0.01733176511849817, 3, 1, 1.605544, 4.816632, 0.4102779191519537, 0.4102779191519537, 1.972640396792, 17.753763571128, 0.007110840528047839, 0.007110840528047842)
    8 = @("sCs", "Sema4d", "Met", "ECs", 3, 1, 1.154206333333333, 3.462619, 0.01628172957288352, 0.01628172957288353, 3, 1, 2.254050666666667, 6.762152, 0.5759961839619929, 0.575996183961993, 2.601639555120889, 23.414755996088, 0.009378214102282038, 0.009378214102282042)
    9 = @("sCs", "Sema4d", "Met", "FAPs", 3, 1, 1.154206333333333, 3.462619, 0.01628172957288352, 0.01628172957288353, 1, 0.3333333333333333, 0.05371366666666667, 0.161141, 0.01372589688605336, 0.01372589688605336, 0.06199665425322223, 0.5579698882790001, 0.0002234813412440049, 0.0002234813412440049)
    10 = @("sCs", "Sema4d", "Met", "sCs", 3, 1, 1.154206333333333, 3.462619, 0.01628172957288352, 0.01628172957288353, 3, 1, 1.605544, 4.816632, 0.4102779191519537, 0.4102779191519537, 1.853129053245334, 16.678161479208, 0.006680034129357479, 0.006680034129357482)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $ws.Cells.Item([int]$r, $i + 1).Value = $rowVals[$i]
    }
}
